$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", copying the same formatting used by the
# other header cells (bold, border, centered) without creating a new style.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New "Save" data column values for rows 2-12.
$saveValues = @(1, 0, 0, 0, 0, 1, 1, 1, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
